$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1
$ws.Range("L2").Value = "stimuli/img_d26ik.png"
$ws.Range("M2").Value = 77.73809523809524
$ws.Range("N2").Value = 60.66666666666666
$ws.Range("O2").Value = 69.20238095238095
$ws.Range("P2").Value = 42
# Row 3
$ws.Range("F3").Value = 2
$ws.Range("L3").Value = "stimuli/img_es7o2.png"
$ws.Range("M3").Value = 52.48571428571429
$ws.Range("N3").Value = 27.54285714285714
$ws.Range("O3").Value = 40.01428571428572
$ws.Range("P3").Value = 35
# Row 4
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = "kitchens"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_6nbgt.png"
$ws.Range("M4").Value = 78.45161290322581
$ws.Range("N4").Value = 57.83870967741935
$ws.Range("O4").Value = 68.14516129032258
$ws.Range("P4").Value = 31
# Row 5
$ws.Range("F5").Value = 4
$ws.Range("H5").Value = "kitchens"
$ws.Range("I5").Value = "target"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_inqod.png"
$ws.Range("M5").Value = 70.84848484848484
$ws.Range("N5").Value = 50.63636363636363
$ws.Range("O5").Value = 60.74242424242424
$ws.Range("P5").Value = 33
$ws.Range("Q5").Value = 5
$ws.Range("R5").Value = 5
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 5
$ws.Range("V5").Value = 5
# Row 6
$ws.Range("F6").Value = 5
$ws.Range("H6").Value = "bedrooms"
$ws.Range("I6").Value = "distractor"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_d9ogj.png"
$ws.Range("M6").Value = 76.86842105263158
$ws.Range("N6").Value = 53.5
$ws.Range("O6").Value = 65.18421052631578
$ws.Range("P6").Value = 38
$ws.Range("Q6").Value = 6
$ws.Range("R6").Value = 6
$ws.Range("S6").Value = 6
$ws.Range("T6").Value = 6
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 6
# Row 7
$ws.Range("F7").Value = 6
$ws.Range("L7").Value = "stimuli/img_5m6x4.png"
$ws.Range("M7").Value = 80.23076923076923
$ws.Range("N7").Value = 58.41025641025641
$ws.Range("O7").Value = 69.32051282051282
$ws.Range("P7").Value = 39
# Row 8
$ws.Range("F8").Value = 7
# Row 9
$ws.Range("F9").Value = 8
$ws.Range("L9").Value = "stimuli/img_fea1z.png"
$ws.Range("M9").Value = 79.45945945945945
$ws.Range("N9").Value = 56.24324324324324
$ws.Range("O9").Value = 67.85135135135135
$ws.Range("P9").Value = 37
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 7
$ws.Range("T9").Value = 7
$ws.Range("U9").Value = 7
$ws.Range("V9").Value = 7
# Row 10
$ws.Range("F10").Value = 9
$ws.Range("H10").Value = "living_rooms"
$ws.Range("L10").Value = "stimuli/img_f63yi.png"
$ws.Range("M10").Value = 85.275
$ws.Range("N10").Value = 68.475
$ws.Range("O10").Value = 76.875
$ws.Range("P10").Value = 40
$ws.Range("Q10").Value = 9
$ws.Range("R10").Value = 9
$ws.Range("S10").Value = 9
$ws.Range("T10").Value = 9
$ws.Range("U10").Value = 8
$ws.Range("V10").Value = 9
# Row 11
$ws.Range("F11").Value = 10
$ws.Range("L11").Value = "stimuli/img_nyv2b.png"
$ws.Range("M11").Value = 11.91176470588235
$ws.Range("N11").Value = 6.852941176470588
$ws.Range("O11").Value = 9.382352941176471
$ws.Range("P11").Value = 34
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 1
$ws.Range("V11").Value = 1
# Row 12
$ws.Range("F12").Value = 11
$ws.Range("H12").Value = "bedrooms"
$ws.Range("I12").Value = "distractor"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_u1rxv.png"
$ws.Range("M12").Value = 75.63636363636364
$ws.Range("N12").Value = 54.27272727272727
$ws.Range("O12").Value = 64.95454545454545
$ws.Range("P12").Value = 44
$ws.Range("Q12").Value = 6
$ws.Range("R12").Value = 6
$ws.Range("S12").Value = 6
$ws.Range("T12").Value = 6
$ws.Range("U12").Value = 6
$ws.Range("V12").Value = 6
# Row 13
$ws.Range("F13").Value = 12
$ws.Range("H13").Value = "kitchens"
$ws.Range("I13").Value = "target"
$ws.Range("K13").Value = "j"
$ws.Range("L13").Value = "stimuli/img_aplao.png"
$ws.Range("M13").Value = 64.0909090909091
$ws.Range("N13").Value = 40.75757575757576
$ws.Range("O13").Value = 52.42424242424242
$ws.Range("P13").Value = 33
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 3
$ws.Range("S13").Value = 3
$ws.Range("T13").Value = 3
$ws.Range("U13").Value = 3
$ws.Range("V13").Value = 3
# Row 14
$ws.Range("F14").Value = 13
$ws.Range("L14").Value = "stimuli/img_0jzz7.png"
$ws.Range("M14").Value = 84.85106382978724
$ws.Range("N14").Value = 68.87234042553192
$ws.Range("O14").Value = 76.86170212765958
$ws.Range("P14").Value = 47
$ws.Range("T14").Value = 9
$ws.Range("V14").Value = 9
# Row 15
$ws.Range("F15").Value = 14
$ws.Range("L15").Value = "stimuli/img_iyxnj.png"
$ws.Range("M15").Value = 75.30555555555556
$ws.Range("N15").Value = 54.33333333333334
$ws.Range("O15").Value = 64.81944444444444
$ws.Range("P15").Value = 36
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = 6
$ws.Range("S15").Value = 6
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = 6
$ws.Range("V15").Value = 6
# Row 16
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = "kitchens"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_t90e2.png"
$ws.Range("M16").Value = 83.0625
$ws.Range("N16").Value = 61.96875
$ws.Range("O16").Value = 72.515625
$ws.Range("P16").Value = 32
$ws.Range("Q16").Value = 9
$ws.Range("R16").Value = 9
$ws.Range("S16").Value = 9
$ws.Range("T16").Value = 9
$ws.Range("U16").Value = 9
$ws.Range("V16").Value = 9
# Row 17
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = "bedrooms"
$ws.Range("I17").Value = "distractor"
$ws.Range("K17").Value = "f"
$ws.Range("L17").Value = "stimuli/img_iqmdm.png"
$ws.Range("M17").Value = 79.38888888888889
$ws.Range("N17").Value = 58.36111111111111
$ws.Range("O17").Value = 68.875
$ws.Range("P17").Value = 36
$ws.Range("Q17").Value = 7
$ws.Range("R17").Value = 7
$ws.Range("S17").Value = 7
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 7
$ws.Range("V17").Value = 7
# Row 18
$ws.Range("F18").Value = 17
$ws.Range("L18").Value = "stimuli/img_yeh72.png"
$ws.Range("M18").Value = 68.66666666666667
$ws.Range("N18").Value = 45.21212121212121
$ws.Range("O18").Value = 56.93939393939394
$ws.Range("P18").Value = 33
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 4
$ws.Range("S18").Value = 4
$ws.Range("T18").Value = 4
$ws.Range("U18").Value = 4
$ws.Range("V18").Value = 4
# Row 19
$ws.Range("F19").Value = 18
$ws.Range("H19").Value = "kitchens"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_p659z.png"
$ws.Range("M19").Value = 84.21621621621621
$ws.Range("N19").Value = 65.37837837837837
$ws.Range("O19").Value = 74.79729729729729
$ws.Range("P19").Value = 37
# Row 20
$ws.Range("F20").Value = 19
$ws.Range("L20").Value = "stimuli/img_30vz5.png"
$ws.Range("M20").Value = 86.21212121212122
$ws.Range("N20").Value = 68.27272727272727
$ws.Range("O20").Value = 77.24242424242425
$ws.Range("P20").Value = 33
$ws.Range("Q20").Value = 10
$ws.Range("R20").Value = 10
$ws.Range("S20").Value = 10
$ws.Range("T20").Value = 10
$ws.Range("U20").Value = 10
$ws.Range("V20").Value = 10
# Row 21
$ws.Range("F21").Value = 20
$ws.Range("H21").Value = "living_rooms"
$ws.Range("I21").Value = "distractor"
$ws.Range("K21").Value = "f"
$ws.Range("L21").Value = "stimuli/img_wgkqa.png"
$ws.Range("M21").Value = 87.25581395348837
$ws.Range("N21").Value = 71.13953488372093
$ws.Range("O21").Value = 79.19767441860465
$ws.Range("P21").Value = 43
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = 9
$ws.Range("U21").Value = 9
$ws.Range("V21").Value = 9
# Row 22
$ws.Range("F22").Value = 21
$ws.Range("L22").Value = "stimuli/img_wyl6z.png"
$ws.Range("M22").Value = 59.8235294117647
$ws.Range("N22").Value = 36.23529411764706
$ws.Range("O22").Value = 48.02941176470588
$ws.Range("P22").Value = 34
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 3
$ws.Range("S22").Value = 3
$ws.Range("T22").Value = 3
$ws.Range("U22").Value = 3
$ws.Range("V22").Value = 3
# Row 23
$ws.Range("F23").Value = 22
$ws.Range("L23").Value = "stimuli/img_yosqb.png"
$ws.Range("M23").Value = 50.88372093023256
$ws.Range("N23").Value = 30.11627906976744
$ws.Range("O23").Value = 40.5
$ws.Range("P23").Value = 43
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 3
$ws.Range("S23").Value = 3
$ws.Range("T23").Value = 3
$ws.Range("U23").Value = 3
$ws.Range("V23").Value = 3
# Row 24
$ws.Range("F24").Value = 23
$ws.Range("H24").Value = "bedrooms"
$ws.Range("I24").Value = "distractor"
$ws.Range("K24").Value = "f"
$ws.Range("L24").Value = "stimuli/img_n9xll.png"
$ws.Range("M24").Value = 77.14285714285714
$ws.Range("N24").Value = 59.21428571428572
$ws.Range("O24").Value = 68.17857142857143
$ws.Range("P24").Value = 42
$ws.Range("Q24").Value = 7
$ws.Range("R24").Value = 7
$ws.Range("S24").Value = 7
$ws.Range("T24").Value = 7
$ws.Range("U24").Value = 7
$ws.Range("V24").Value = 7
# Row 25
$ws.Range("F25").Value = 24
$ws.Range("H25").Value = "living_rooms"
$ws.Range("I25").Value = "distractor"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_pna7l.png"
$ws.Range("M25").Value = 85.53333333333333
$ws.Range("N25").Value = 67.97777777777777
$ws.Range("O25").Value = 76.75555555555556
$ws.Range("P25").Value = 45
$ws.Range("Q25").Value = 9
$ws.Range("R25").Value = 9
$ws.Range("S25").Value = 9
$ws.Range("T25").Value = 9
$ws.Range("U25").Value = 9
$ws.Range("V25").Value = 9
# Row 26
$ws.Range("F26").Value = 25
$ws.Range("L26").Value = "stimuli/img_ye5sl.png"
$ws.Range("M26").Value = 53.2258064516129
$ws.Range("N26").Value = 34.45161290322581
$ws.Range("O26").Value = 43.83870967741936
$ws.Range("P26").Value = 31
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 2
$ws.Range("T26").Value = 2
$ws.Range("U26").Value = 2
$ws.Range("V26").Value = 2
# Row 27
$ws.Range("F27").Value = 26
$ws.Range("H27").Value = "kitchens"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_7wul8.png"
$ws.Range("M27").Value = 43.03030303030303
$ws.Range("N27").Value = 25.54545454545455
$ws.Range("O27").Value = 34.28787878787879
$ws.Range("P27").Value = 33
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("U27").Value = 1
$ws.Range("V27").Value = 1
# Row 28
$ws.Range("F28").Value = 27
$ws.Range("H28").Value = "kitchens"
$ws.Range("I28").Value = "target"
$ws.Range("K28").Value = "j"
$ws.Range("L28").Value = "stimuli/img_d8xbu.png"
$ws.Range("M28").Value = 91.36363636363636
$ws.Range("N28").Value = 73.18181818181819
$ws.Range("O28").Value = 82.27272727272728
$ws.Range("P28").Value = 33
$ws.Range("T28").Value = 10
$ws.Range("U28").Value = 10
$ws.Range("V28").Value = 10
# Row 29
$ws.Range("F29").Value = 28
$ws.Range("L29").Value = "stimuli/img_njmgp.png"
$ws.Range("M29").Value = 80.48148148148148
$ws.Range("N29").Value = 58.4074074074074
$ws.Range("O29").Value = 69.44444444444444
$ws.Range("P29").Value = 27
$ws.Range("Q29").Value = 8
$ws.Range("R29").Value = 8
$ws.Range("S29").Value = 8
$ws.Range("T29").Value = 8
$ws.Range("U29").Value = 8
$ws.Range("V29").Value = 8
# Row 30
$ws.Range("F30").Value = 29
$ws.Range("H30").Value = "bedrooms"
$ws.Range("L30").Value = "stimuli/img_8dacu.png"
$ws.Range("M30").Value = 76.38461538461539
$ws.Range("N30").Value = 53.64102564102564
$ws.Range("O30").Value = 65.01282051282051
$ws.Range("P30").Value = 39
$ws.Range("Q30").Value = 6
$ws.Range("R30").Value = 6
$ws.Range("S30").Value = 6
$ws.Range("T30").Value = 6
$ws.Range("U30").Value = 6
$ws.Range("V30").Value = 6
# Row 31
$ws.Range("F31").Value = 30
$ws.Range("L31").Value = "stimuli/img_ac0ey.png"
$ws.Range("M31").Value = 86.62222222222222
$ws.Range("N31").Value = 70.02222222222223
$ws.Range("O31").Value = 78.32222222222222
$ws.Range("P31").Value = 45
$ws.Range("U31").Value = 9
# Row 32
$ws.Range("F32").Value = 31
$ws.Range("H32").Value = "bedrooms"
$ws.Range("L32").Value = "stimuli/img_ybbmx.png"
$ws.Range("M32").Value = 55.24324324324324
$ws.Range("N32").Value = 36.75675675675676
$ws.Range("O32").Value = 46
$ws.Range("P32").Value = 37
$ws.Range("Q32").Value = 3
$ws.Range("R32").Value = 3
$ws.Range("S32").Value = 3
$ws.Range("T32").Value = 3
$ws.Range("U32").Value = 3
$ws.Range("V32").Value = 3
# Row 33
$ws.Range("F33").Value = 32
$ws.Range("L33").Value = "stimuli/img_cv6mf.png"
$ws.Range("M33").Value = 66.8
$ws.Range("N33").Value = 42.08
$ws.Range("O33").Value = 54.44
$ws.Range("P33").Value = 25
# Row 34
$ws.Range("F34").Value = 33
$ws.Range("H34").Value = "living_rooms"
$ws.Range("L34").Value = "stimuli/img_6zz63.png"
$ws.Range("M34").Value = 87.66666666666667
$ws.Range("N34").Value = 70.6
$ws.Range("O34").Value = 79.13333333333333
$ws.Range("P34").Value = 45
$ws.Range("Q34").Value = 9
$ws.Range("R34").Value = 10
$ws.Range("S34").Value = 10
$ws.Range("T34").Value = 9
$ws.Range("U34").Value = 9
$ws.Range("V34").Value = 9
# Row 35
$ws.Range("F35").Value = 34
# Row 36
$ws.Range("F36").Value = 35
$ws.Range("H36").Value = "bedrooms"
$ws.Range("I36").Value = "distractor"
$ws.Range("K36").Value = "f"
$ws.Range("L36").Value = "stimuli/img_dkqas.png"
$ws.Range("M36").Value = 78.57894736842105
$ws.Range("N36").Value = 57.71052631578947
$ws.Range("O36").Value = 68.14473684210526
$ws.Range("P36").Value = 38
$ws.Range("Q36").Value = 7
$ws.Range("R36").Value = 7
$ws.Range("S36").Value = 7
$ws.Range("T36").Value = 7
$ws.Range("U36").Value = 7
$ws.Range("V36").Value = 7
# Row 37
$ws.Range("F37").Value = 36
$ws.Range("L37").Value = "stimuli/img_cnyac.png"
$ws.Range("M37").Value = 69.1470588235294
$ws.Range("N37").Value = 47.8235294117647
$ws.Range("O37").Value = 58.48529411764706
$ws.Range("P37").Value = 34
$ws.Range("Q37").Value = 5
$ws.Range("R37").Value = 5
$ws.Range("S37").Value = 5
$ws.Range("T37").Value = 5
$ws.Range("U37").Value = 5
$ws.Range("V37").Value = 5
# Row 38
$ws.Range("F38").Value = 37
$ws.Range("H38").Value = "kitchens"
$ws.Range("I38").Value = "target"
$ws.Range("K38").Value = "j"
$ws.Range("L38").Value = "stimuli/img_eatdk.png"
$ws.Range("M38").Value = 81.40625
$ws.Range("N38").Value = 61.375
$ws.Range("O38").Value = 71.390625
$ws.Range("P38").Value = 32
$ws.Range("Q38").Value = 8
$ws.Range("R38").Value = 8
$ws.Range("S38").Value = 8
$ws.Range("T38").Value = 8
$ws.Range("U38").Value = 8
$ws.Range("V38").Value = 8
# Row 39
$ws.Range("F39").Value = 38
$ws.Range("H39").Value = "kitchens"
$ws.Range("I39").Value = "target"
$ws.Range("K39").Value = "j"
$ws.Range("L39").Value = "stimuli/img_p3hpc.png"
$ws.Range("M39").Value = 72.83333333333333
$ws.Range("N39").Value = 52.22222222222222
$ws.Range("O39").Value = 62.52777777777777
$ws.Range("P39").Value = 36
$ws.Range("Q39").Value = 6
$ws.Range("R39").Value = 6
$ws.Range("S39").Value = 6
$ws.Range("T39").Value = 6
$ws.Range("U39").Value = 6
$ws.Range("V39").Value = 6
# Row 40
$ws.Range("F40").Value = 39
$ws.Range("H40").Value = "living_rooms"
$ws.Range("I40").Value = "distractor"
$ws.Range("K40").Value = "f"
$ws.Range("L40").Value = "stimuli/img_amsgw.png"
$ws.Range("M40").Value = 86.08510638297872
$ws.Range("N40").Value = 65.95744680851064
$ws.Range("O40").Value = 76.02127659574468
$ws.Range("P40").Value = 47
$ws.Range("Q40").Value = 9
$ws.Range("R40").Value = 9
$ws.Range("S40").Value = 9
$ws.Range("T40").Value = 8
$ws.Range("U40").Value = 9
$ws.Range("V40").Value = 8
# Row 41
$ws.Range("F41").Value = 40
$ws.Range("H41").Value = "living_rooms"
$ws.Range("I41").Value = "distractor"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_95hiv.png"
$ws.Range("M41").Value = 84.04545454545455
$ws.Range("N41").Value = 67.31818181818181
$ws.Range("O41").Value = 75.68181818181819
$ws.Range("P41").Value = 44
$ws.Range("Q41").Value = 9
$ws.Range("R41").Value = 9
$ws.Range("S41").Value = 9
$ws.Range("T41").Value = 8
$ws.Range("U41").Value = 8
$ws.Range("V41").Value = 8
